$wb = $excel.ActiveWorkbook

# Add new worksheet "ID management" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "ID management"

# Populate the data
$ws2.Range("A1").Value = "When"
$ws2.Range("B1").Value = "then id ="

$ws2.Range("A2").Value = "new Ingredient"
$ws2.Range("B2").Value = "FoodItem.id"

$ws2.Range("A3").Value = "new ComposedFoodItem"
$ws2.Range("B3").Value = "new"

$ws2.Range("A4").Value = "FoodItem from FoodItemVM"
$ws2.Range("B4").Value = "FoodItemVM.id"

$ws2.Range("A5").Value = "FoodItem from ComposedFoodItemVM"
$ws2.Range("B5").Value = "new"

# Column A width (best-fit similar to source; engine quantizes to nearest
# pixel-based width, so this lands as close as possible to 32.6640625)
$ws2.Columns.Item(1).ColumnWidth = 31.83

# Selection on the new sheet
$ws2.Range("A6").Select() | Out-Null
